$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "11.00", "0.999").
# Force text storage via a temporary "@" (Text) number format, then restore
# the cell to the workbook default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.838.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.333.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.330.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.183'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.10%  '
$ws.Range("E11").Value = '  +2.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.61%  '
$ws.Range("E13").Value = '  +3.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '698.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.876.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.835.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.336.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.895'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.53%  '
$ws.Range("E27").Value = '  +2.68%  '
$ws.Range("E28").Value = '  +6.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.99'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("E30").Value = '  +3.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '567.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +3.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.687.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.88%  '
$ws.Range("E40").Value = '  +5.29%  '
$ws.Range("E41").Value = '  +3.74%  '
$ws.Range("E42").Value = '  +7.83%  '
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("E44").Value = '  +4.46%  '
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("E46").Value = '  +2.78%  '
$ws.Range("E47").Value = '  +6.36%  '
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.62%  '
